$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header A1: "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Column A data values (row 2-14): Gen counts -> MaxFES fractions
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# 3. Column L (Run 10) previously-blank cells now carry values
$ws.Range("L3").Value = 163306175.228072
$ws.Range("L4").Value = 73959293.55376348
$ws.Range("L5").Value = 13422.10640269
$ws.Range("L6").Value = 2214.97578219
$ws.Range("L7").Value = 798.0581635
$ws.Range("L8").Value = 64.81965723
$ws.Range("L9").Value = 0.57616754
$ws.Range("L10").Value = 0.06929021
$ws.Range("L11").Value = 0.00045578
$ws.Range("L12").Value = 0.0000364
$ws.Range("L13").Value = 0.00000508
$ws.Range("L14").Value = 0.00000001

# 4. Drop the "Run 50" column entirely. This shifts the trailing "Mean"
#    column (formerly BA) one column left into AZ, matching the diff
#    (dimension A1:BA14 -> A1:AZ14, spans 1:53 -> 1:52).
$ws.Columns("AZ").Delete()

# 5. The Mean column (now AZ) is recalculated to reflect the removed
#    "Run 50" run and the newly-populated "Run 10" (L) column.
$ws.Range("AZ2").Value = 440361151.9764096
$ws.Range("AZ3").Value = 178745815.7647588
$ws.Range("AZ4").Value = 17918068.06625672
$ws.Range("AZ5").Value = 2399555.51343886
$ws.Range("AZ6").Value = 2219834.05472267
$ws.Range("AZ7").Value = 2219795.99612235
$ws.Range("AZ8").Value = 2219781.33133539
$ws.Range("AZ9").Value = 2219780.04646559
$ws.Range("AZ10").Value = 2219780.03632805
$ws.Range("AZ11").Value = 2219780.03495136
$ws.Range("AZ12").Value = 2219780.03494297
$ws.Range("AZ13").Value = 2219780.03494234
$ws.Range("AZ14").Value = 2219780.03494224
